$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.936.37"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +3.52%  '
$ws.Range("D3").Value = "'2.260.93"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.59%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = "'252.30"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.18%  '
$ws.Range("D6").Value = "'0.639"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.13%  '
$ws.Range("D7").Value = "'71.24"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.14%  '
$ws.Range("D8").Value = "'0.670"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +14.86%  '
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").Value = "'39.50"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +4.57%  '
$ws.Range("D11").Value = "'0.0975"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.81%  '
$ws.Range("D12").Value = "'59.60"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.64%  '
$ws.Range("D13").Value = "'7.65"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +6.66%  '
$ws.Range("D14").Value = "'0.104"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.43%  '
$ws.Range("D15").Value = "'2.597.74"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.58%  '
$ws.Range("D16").Value = "'0.890"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.01%  '
$ws.Range("D17").Value = "'14.86"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.79%  '
$ws.Range("D18").Value = "'2.264.62"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.55%  '
$ws.Range("D19").Value = "'42.847.74"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.28%  '
$ws.Range("D20").Value = "'0.0₃0984"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.11%  '
$ws.Range("D21").Value = "'6.29"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.94%  '
$ws.Range("D22").Value = "'73.20"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.85%  '
$ws.Range("D23").Value = "'237.46"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.31%  '
$ws.Range("E24").Value = '  +3.67%  '
$ws.Range("D25").Value = "'3.94"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.74%  '
$ws.Range("D26").Value = "'11.79"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.27%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").Value = "'2.45"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.54%  '
$ws.Range("E29").Value = '  -1.35%  '
$ws.Range("E30").Value = '  +2.27%  '
$ws.Range("D31").Value = "'167.94"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("D32").Value = "'21.24"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.13%  '
$ws.Range("D33").Value = "'6.27"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +13.90%  '
$ws.Range("E34").Value = '  +6.75%  '
$ws.Range("D35").Value = "'0.0773"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.58%  '
$ws.Range("E36").Value = '  +2.50%  '
$ws.Range("D37").Value = "'29.06"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +13.34%  '
$ws.Range("D38").Value = "'4.72"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.62%  '
$ws.Range("D39").Value = "'4.13"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.05%  '
$ws.Range("E40").Value = '  +6.86%  '
$ws.Range("E41").Value = '  +3.42%  '
$ws.Range("D42").Value = "'5.85"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.18%  '
$ws.Range("D43").Value = "'12.17"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").Value = "'64.59"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.44%  '
$ws.Range("D45").Value = "'5.03"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.41%  '
$ws.Range("D46").Value = "'0.203"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.92%  '
$ws.Range("D47").Value = "'8.93"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.61%  '
$ws.Range("D48").Value = "'0.104"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.51%  '
$ws.Range("E49").Value = '  -5.07%  '
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("D51").Value = "'1.19"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.91%  '